$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-FormattedCell($targetAddr, $sourceAddr, $value) {
    $ws.Range($sourceAddr).Copy()
    $ws.Range($targetAddr).PasteSpecial(-4104)  # xlPasteAll
    $excel.CutCopyMode = $false
    $ws.Range($targetAddr).Value2 = $value
}

# Row 18: SSD02, 2020-05-13 13:00-15:00
Set-FormattedCell "A18" "A14" "SSD02"
Set-FormattedCell "C18" "C14" 43964
Set-FormattedCell "D18" "D14" 0.54166666666666663
Set-FormattedCell "E18" "E14" 0.625

# Row 19: småting, 2020-05-13 15:00-16:00
Set-FormattedCell "A19" "A14" "småting"
Set-FormattedCell "C19" "C14" 43964
Set-FormattedCell "D19" "D14" 0.625
Set-FormattedCell "E19" "E14" 0.66666666666666663

# Row 20: Hjælp af Anders, 2020-05-13 16:00-17:40
Set-FormattedCell "A20" "A14" "Hjælp af Anders"
Set-FormattedCell "C20" "C14" 43964
Set-FormattedCell "D20" "D14" 0.66666666666666663
Set-FormattedCell "E20" "E14" 0.73611111111111116

# Update view state: scroll to A7 and select A22
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("A22").Select()
